$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
$shp = $np.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "Hallo Zusammen `rIch lade Sie zu meiner HTML-Präsentation „Lieblingsgerichte meiner Familie“ herzlich ein."
Write-Host "Done: $($shp.TextFrame.TextRange.Text)"
